# Apply trade #70 close-out update across the workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.57   # Current Capital
$wsSummary.Range("B4").Value = 0.58      # Total P&L $
$wsSummary.Range("B5").Value = 0.17      # Total P&L %
$wsSummary.Range("B6").Value = 70        # Total Trades
$wsSummary.Range("B8").Value = 28        # Losing Trades
$wsSummary.Range("B9").Value = 41.43     # Win Rate %

# --- Strategy Status sheet -------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.57     # Capital
$wsStatus.Range("D4").Value = 70         # Trades
$wsStatus.Range("E4").Value = 0.58       # P&L $
$wsStatus.Range("F4").Value = 0.57       # P&L %
$wsStatus.Range("G4").Value = 41.43      # Win Rate %

# --- Helper: append the new trade #70 row to a trades sheet ----------
function Add-Trade70Row($ws) {
    $ws.Range("A71").Value = 70

    # Force literal text so "2026-02-17" isn't auto-parsed into a date serial.
    $ws.Range("B71").NumberFormat = "@"
    $ws.Range("B71").Value = "2026-02-17"

    $ws.Range("C71").Value = "08:57:21"
    $ws.Range("D71").Value = "MarketMaking"
    $ws.Range("E71").Value = "UP"
    $ws.Range("F71").Value = 0.14
    $ws.Range("G71").Value = 0.1
    $ws.Range("H71").Value = "CLOSED"
    $ws.Range("I71").Value = -28.5714
    $ws.Range("J71").Value = -0.04
    $ws.Range("K71").Value = 100.57
    $ws.Range("L71").Value = 0
    $ws.Range("M71").Value = 0
    $ws.Range("N71").Value = 0.6
    $ws.Range("O71").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P71").Value = "early_exit"
    $ws.Range("Q71").Value = 0.1
}

# --- All Trades sheet --------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade70Row $wsAllTrades

# --- MarketMaking sheet -------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade70Row $wsMarketMaking
